# Scheduled-runner price/profit refresh across the Leve Profits sheets.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H,I,J,K,L,M,N) for the rows whose market prices moved since the
# last run. A few rows gain a previously-blank LeveProfitHQ (N) value, and a
# few rows lose a previously-populated LeveProfit value (cleared, because
# that recipe no longer has an HQ/NQ price to diff against).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 434
$ws.Range("I20").Value = 434
$ws.Range("K20").Value = 434
$ws.Range("M20").Value = -204
$ws.Range("H33").Value = 820.5
$ws.Range("I33").Value = 94.166664
$ws.Range("K33").Value = 94.166664
$ws.Range("M33").Value = 134.833336
$ws.Range("H35").Value = 434
$ws.Range("I35").Value = 434
$ws.Range("K35").Value = 434
$ws.Range("M35").Value = -55
$ws.Range("H48").Value = 9507
$ws.Range("J48").Value = 9507
$ws.Range("L48").Value = 28521
$ws.Range("N48").Value = -29105
$ws.Range("H53").Value = 120.92308
$ws.Range("I53").Value = 109.333336
$ws.Range("K53").Value = 109.333336
$ws.Range("M53").Value = 527.666664
$ws.Range("H56").Value = 9507
$ws.Range("J56").Value = 9507
$ws.Range("L56").Value = 28521
$ws.Range("N56").Value = -29589
$ws.Range("H74").Value = 4497
$ws.Range("J74").Value = 5000
$ws.Range("L74").Value = 5000
$ws.Range("N74").Value = -6872
$ws.Range("H77").Value = 4497
$ws.Range("J77").Value = 5000
$ws.Range("L77").Value = 25000
$ws.Range("N77").Value = -34360
$ws.Range("H125").Value = 3354.4
$ws.Range("I125").Value = 2540.4285
$ws.Range("K125").Value = 22863.8565
$ws.Range("M125").Value = -20403.8565
$ws.Range("H132").Value = 2519.9412
$ws.Range("I132").Value = 2519.9412
$ws.Range("K132").Value = 7559.823600000001
$ws.Range("M132").Value = -5029.823600000001
$ws.Range("H141").Value = 1561.5385
$ws.Range("I141").Value = 1561.5385
$ws.Range("K141").Value = 4684.6155
$ws.Range("M141").Value = 495.3845000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2875
$ws.Range("I45").Value = 3250
$ws.Range("K45").Value = 3250
$ws.Range("M45").Value = -2873
$ws.Range("H110").Value = 524.6667
$ws.Range("I110").Value = 419.5
$ws.Range("J110").Value = 735
$ws.Range("K110").Value = 419.5
$ws.Range("L110").Value = 735
$ws.Range("M110").Value = 1625.5
$ws.Range("N110").Value = -4825
$ws.Range("H122").Value = 2497.75
$ws.Range("J122").Value = 4999
$ws.Range("L122").Value = 14997
$ws.Range("N122").Value = -19897

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 28500
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H107").Value = 2327
$ws.Range("I107").Value = 2269.3333
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 2269.3333
$ws.Range("L107").Value = 2500
$ws.Range("M107").Value = -349.3332999999998
$ws.Range("N107").Value = -6340
$ws.Range("H134").Value = 7256.8
$ws.Range("I134").Value = 6954.579
$ws.Range("J134").Value = 12999
$ws.Range("K134").Value = 20863.737
$ws.Range("L134").Value = 38997
$ws.Range("M134").Value = -18328.737
$ws.Range("N134").Value = -44067

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5114.6665
$ws.Range("I99").Value = 3566.5
$ws.Range("K99").Value = 3566.5
$ws.Range("M99").Value = -2068.5
$ws.Range("H126").Value = 5114.6665
$ws.Range("I126").Value = 3566.5
$ws.Range("K126").Value = 10699.5
$ws.Range("M126").Value = -8229.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 470.83334
$ws.Range("I9").Value = 441.66666
$ws.Range("K9").Value = 1324.99998
$ws.Range("M9").Value = -1100.99998
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H12").Value = 116.875
$ws.Range("H23").Value = 352.0909
$ws.Range("I23").Value = 477
$ws.Range("J23").Value = 248
$ws.Range("K23").Value = 1431
$ws.Range("L23").Value = 744
$ws.Range("M23").Value = -1196
$ws.Range("N23").Value = -1214
$ws.Range("H38").Value = 179.76923
$ws.Range("I38").Value = 160.5
$ws.Range("J38").Value = 210.6
$ws.Range("K38").Value = 481.5
$ws.Range("L38").Value = 631.8
$ws.Range("M38").Value = -134.5
$ws.Range("N38").Value = -1325.8
$ws.Range("H98").Value = 99
$ws.Range("I98").Value = 99
$ws.Range("K98").Value = 297
$ws.Range("M98").Value = 1201
$ws.Range("H131").Value = 732.5
$ws.Range("J131").Value = 700
$ws.Range("L131").Value = 2100
$ws.Range("N131").Value = -12180
$ws.Range("H132").Value = 3008.4
$ws.Range("I132").Value = 2009.3334
$ws.Range("K132").Value = 18084.0006
$ws.Range("M132").Value = -15554.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 2934
$ws.Range("I9").Value = 512
$ws.Range("J9").Value = 10200
$ws.Range("K9").Value = 512
$ws.Range("L9").Value = 10200
$ws.Range("M9").Value = -342
$ws.Range("N9").Value = -10540
$ws.Range("H100").Value = 37349.75
$ws.Range("J100").Value = 37349.75
$ws.Range("L100").Value = 37349.75
$ws.Range("N100").Value = -39513.75
$ws.Range("H102").Value = 2331
$ws.Range("I102").Value = 2998
$ws.Range("K102").Value = 2998
$ws.Range("M102").Value = -1376
$ws.Range("H122").Value = 2439.5
$ws.Range("I122").Value = 2439.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7318.5
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -4868.5
$ws.Range("H126").Value = 18933
$ws.Range("I126").Value = 18933
$ws.Range("K126").Value = 56799
$ws.Range("M126").Value = -54329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 38242.75
$ws.Range("J7").Value = 27250
$ws.Range("L7").Value = 27250
$ws.Range("N7").Value = -27474
$ws.Range("H16").Value = 430.5
$ws.Range("I16").Value = 430.5
$ws.Range("K16").Value = 430.5
$ws.Range("M16").Value = -260.5
$ws.Range("H120").Value = 59997
$ws.Range("J120").Value = 59997
$ws.Range("L120").Value = 59997
$ws.Range("N120").Value = -69673
$ws.Range("H126").Value = 38242.75
$ws.Range("J126").Value = 27250
$ws.Range("L126").Value = 81750
$ws.Range("N126").Value = -86690
